# Refresh the Price (D) and Volume(1h) (E) columns with the latest scrape
# values, as produced by the scheduled GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.241.32'
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").Value = '1.898.10'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = "'307.88"
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("D7").Value = "'0.5211"
$ws.Range("E7").Value = '  +0.73%  '
$ws.Range("D8").Value = "'0.3773"
$ws.Range("E8").Value = '  +0.50%  '
$ws.Range("D9").Value = "'0.07283"
$ws.Range("E9").Value = '  +1.04%  '
$ws.Range("D10").Value = "'21.21"
$ws.Range("E10").Value = '  +0.51%  '
$ws.Range("D11").Value = "'0.9017"
$ws.Range("E11").Value = '  +0.56%  '
$ws.Range("D12").Value = "'0.08176"
$ws.Range("E12").Value = '  +6.78%  '
$ws.Range("D13").Value = "'96.67"
$ws.Range("E13").Value = '  +2.53%  '
$ws.Range("D14").Value = '1.903.04'
$ws.Range("E14").Value = '  +0.59%  '
$ws.Range("D15").Value = "'5.295"
$ws.Range("E15").Value = '  +1.26%  '
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = '  +0.21%  '
$ws.Range("D17").Value = "'0.000008611"
$ws.Range("E17").Value = '  +1.13%  '
$ws.Range("E18").Value = '  +1.27%  '
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("D20").Value = '27.276.54'
$ws.Range("E20").Value = '  +0.52%  '
$ws.Range("D21").Value = "'5.095"
$ws.Range("E21").Value = '  +0.80%  '
$ws.Range("D22").Value = "'10.71"
$ws.Range("D23").Value = "'6.412"
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  +1.04%  '
$ws.Range("D26").Value = "'18.25"
$ws.Range("E26").Value = '  +1.12%  '
$ws.Range("E27").Value = '  +0.77%  '
$ws.Range("D28").Value = "'115.43"
$ws.Range("E28").Value = '  +0.99%  '
$ws.Range("D29").Value = "'4.827"
$ws.Range("E29").Value = '  +1.11%  '
$ws.Range("D30").Value = "'4.915"
$ws.Range("E30").Value = '  -1.02%  '
$ws.Range("D31").Value = "'0.09238"
$ws.Range("E31").Value = '  +0.51%  '
$ws.Range("D32").Value = "'0.05057"
$ws.Range("E32").Value = '  +0.27%  '
$ws.Range("E33").Value = '  +2.91%  '
$ws.Range("D34").Value = "'1.230"
$ws.Range("E34").Value = '  -0.50%  '
$ws.Range("D35").Value = "'3.431"
$ws.Range("E35").Value = '  +4.66%  '
$ws.Range("D36").Value = "'2.964"
$ws.Range("E36").Value = '  -0.52%  '
$ws.Range("D37").Value = "'2.589"
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("D38").Value = "'0.5682"
$ws.Range("E38").Value = '  +1.50%  '
$ws.Range("D39").Value = "'0.01996"
$ws.Range("E39").Value = '  +0.40%  '
$ws.Range("D40").Value = "'1.074"
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("D41").Value = "'8.966"
$ws.Range("E41").Value = '  -0.21%  '
$ws.Range("D42").Value = "'6.568"
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("D43").Value = "'115.33"
$ws.Range("E43").Value = '  -3.32%  '
$ws.Range("D44").Value = "'0.1517"
$ws.Range("E44").Value = '  +0.32%  '
$ws.Range("D45").Value = "'0.4882"
$ws.Range("E45").Value = '  +1.32%  '
$ws.Range("E46").Value = '  +0.26%  '
$ws.Range("D47").Value = "'10.05"
$ws.Range("E47").Value = '  -0.72%  '
$ws.Range("D48").Value = "'1.622"
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("D49").Value = "'38.23"
$ws.Range("E49").Value = '  +2.31%  '
$ws.Range("D50").Value = "'63.62"
$ws.Range("E50").Value = '  -0.31%  '
$ws.Range("D51").Value = "'0.05943"
$ws.Range("E51").Value = '  +0.32%  '

# Drop the quote-prefix style Excel just applied to the text-forced cells
# above so they keep the sheet's original (unstyled) formatting.
$ws.Range("D4").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
